$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Refcode" value in J2 (test name replaced with tester name)
$ws.Range("J2").Value = "tester23"

# Update the "Email Name" value in A2 (EOM code bumped from Feb to Mar)
$ws.Range("A2").Value = "3.25.21.EOM2"

# Move the active selection from E2 to C3
$ws.Range("C3").Select()
